$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Map of cell address -> new value. Values are kept as text so that
# numeric-looking strings (prices) do not get silently reinterpreted
# as numbers by Excel (the source data stores these as text).
$updates = [ordered]@{
    'D2' = '245.52'
    'D3' = '23.95'
    'D4' = '5.164'
    'D5' = '0.05736'
    'D6' = '6.488'
    'D7' = '3.164'
    'D8' = '0.8158'
    'D9' = '0.8494'
    'D10' = '0.1369'
    'D11' = '0.06952'
    'D12' = '0.03179'
    'D14' = '0.09340'
    'D15' = '3.821'
    'D16' = '0.001518'
    'D17' = '0.04689'
    'D18' = '0.0005960'
    'E18' = '17OneONEWorstin24h'
    'D19' = '0.006193'
    'D20' = '0.001242'
    'D21' = '0.004830'
    'D22' = '0.00008499'
    'D23' = '3.541'
    'D24' = '2.158'
    'D25' = '0.3199'
    'D26' = '0.1336'
    'D27' = '0.0002327'
    'D40' = '0.03697'
    'D41' = '0.006429'
    'E41' = '40KickTokenKICK'
    'D42' = '0.1056'
    'D43' = '0.002259'
    'D44' = '0.007798'
    'E44' = '43LocalTradersLCTBestin24h'
    'D45' = '0.00005463'
    'D47' = '0.3993'
    'D48' = '0.002557'
    'E48' = '47BOLOBOLO'
    'D49' = '0.00002096'
    'D50' = '0.0001996'
}

foreach ($addr in $updates.Keys) {
    $cell = $ws.Range($addr)
    $originalStyle = $cell.Style
    # Force text interpretation so numeric-looking values (e.g. "245.52")
    # are stored as text, matching the original cell type, then restore
    # the cell style so no formatting is changed by this trick.
    $cell.NumberFormat = "@"
    $cell.Value = $updates[$addr]
    $cell.Style = $originalStyle
}
